{"js": "// Insert the missing \"term\u00e9kek \u00f6sszehasonl\u00edt\u00e1sa\" clause into the\n// sentence about personalized recommendations, turning:\n//   \"...licit\u00e1lt felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra.\"\n// into:\n//   \"...licit\u00e1lt felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra, valamint term\u00e9kek \u00f6sszehasonl\u00edt\u00e1s\u00e1ra\n//    is lehet\u0151s\u00e9get kell ny\u00fajtania.\"\nconst body = context.document.body;\n\n// Narrow, unique anchor right at the end of the first sentence (the one\n// that ends \"...licit\u00e1lt felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra.\"). Replacing just this tail\n// keeps every other run in the paragraph untouched.\nconst oldTail = \"felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra.\";\nconst newTail = \"felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra, valamint term\u00e9kek \u00f6sszehasonl\u00edt\u00e1s\u00e1ra is lehet\u0151s\u00e9get kell ny\u00fajtania.\";\n\nconst results = body.search(oldTail, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\n// Use the last match in case the (short) tail phrase happened to recur\n// earlier in the document; here it is expected to be unique.\nconst target = results.items[results.items.length - 1];\ntarget.insertText(newTail, \"Replace\");\nawait context.sync();\n", "ps1": "# Kimaradt a term\u00e9k \u00f6sszehasonl\u00edt\u00e1s...\n# Extend the sentence \"...licit\u00e1lt felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra.\" with the missing\n# clause about product comparison, turning it into:\n# \"...licit\u00e1lt felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra, valamint term\u00e9kek \u00f6sszehasonl\u00edt\u00e1s\u00e1ra\n#  is lehet\u0151s\u00e9get kell ny\u00fajtania.\"\n\n$d = $word.ActiveDocument\n\n$findText    = \"felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra.\"\n$replaceText = \"felhaszn\u00e1l\u00f3k sz\u00e1m\u00e1ra, valamint term\u00e9kek \u00f6sszehasonl\u00edt\u00e1s\u00e1ra is lehet\u0151s\u00e9get kell ny\u00fajtania.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$wdFindWrapContinue = 1\n$wdReplaceAll = 2\n\n$found = $find.Execute(\n    $findText,    # FindText\n    $true,        # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    $wdFindWrapContinue,  # Wrap\n    $false,       # Format\n    $replaceText, # ReplaceWith\n    $wdReplaceAll # Replace\n)\n\nif (-not $found) {\n    throw \"Target sentence end ('$findText') not found in document.\"\n}\n"}
